$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 96, pushing existing rows 96-168 down to 97-169.
$ws.Rows(96).Insert()

# Populate the newly inserted row 96 with the new data record.
$ws.Range("A96").Value = 1
$ws.Range("B96").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C96").Value = "Arica y Parinacota"
$ws.Range("D96").Value = 45086
$ws.Range("E96").Value = 15
$ws.Range("F96").Value = 100112042
$ws.Range("G96").Value = "Locoto"
$ws.Range("H96").Value = "Sin especificar"
$ws.Range("I96").Value = "Segunda"
$ws.Range("J96").Value = 180
$ws.Range("K96").Value = 48000
$ws.Range("L96").Value = 50000
$ws.Range("M96").Value = 48889
$ws.Range("N96").Value = "`$/caja 20 kilos"
$ws.Range("O96").Value = "Región de Arica y Parinacota"
$ws.Range("P96").Value = 2444
$ws.Range("Q96").Value = 20
$ws.Range("R96").Value = "Hortaliza"
